$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two e-mail values that changed.
$ws.Range("B14").Value = "adssfsdfas"
$ws.Range("B17").Value = "dsfafs"

# Move the active selection to match the saved view state.
$ws.Range("B17").Select()
